$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8, 10, 11: H column text fix (remove space after "green")
$ws.Range("H8").Value = "green(raw) petroleum coke (in bulk)"
$ws.Range("H10").Value = "green(raw) petroleum coke (in bulk)"
$ws.Range("H11").Value = "green(raw) petroleum coke (in bulk)"

# Row 12: T/U/V numeric updates
$ws.Range("T12").Value = 2.2023
$ws.Range("U12").Value = 2862.9648
$ws.Range("V12").Value = 3.4846

# Row 14: supplier name change
$ws.Range("O14").Value = "sinopec usa"

# Row 15: T/U updates
$ws.Range("T15").Value = 6.3818
$ws.Range("U15").Value = 6381553.8694

# Row 16: T/U updates
$ws.Range("T16").Value = 6.3818
$ws.Range("U16").Value = 12763107.4665

# Row 17: T/U updates
$ws.Range("T17").Value = 6.3818
$ws.Range("U17").Value = 15953883.9929

# Row 18: T/U updates
$ws.Range("T18").Value = 6.3818
$ws.Range("U18").Value = 1749471.0113

# Row 19: T/U updates
$ws.Range("T19").Value = 6.3818
$ws.Range("U19").Value = 3190776.7986

# Row 31: T/U/V numeric updates
$ws.Range("T31").Value = 2.1837
$ws.Range("U31").Value = 3493.902
$ws.Range("V31").Value = 3.4846

# Row 34: P column supplier name change, T/U/V updates
$ws.Range("P34").Value = "brakes india private limited"
$ws.Range("T34").Value = 84.9866
$ws.Range("U34").Value = 849841.9051
$ws.Range("V34").Value = 0.8075

# Row 35: T/U/V updates
$ws.Range("T35").Value = 83.0891
$ws.Range("U35").Value = 332356.4061
$ws.Range("V35").Value = 0.8075

# Row 36: T/U/V updates
$ws.Range("T36").Value = 83.0891
$ws.Range("U36").Value = 332356.4061
$ws.Range("V36").Value = 0.8075

# Row 37: T/U/V updates
$ws.Range("T37").Value = 671.8729
$ws.Range("U37").Value = 134373.5613
$ws.Range("V37").Value = 8.228199999999999

# Row 38: T/U/V updates
$ws.Range("T38").Value = 79928.6032
$ws.Range("U38").Value = 1918286.4776
$ws.Range("V38").Value = 984.9489

# Row 40: H column item description change
$ws.Range("H40").Value = "calcined petroleum coke"

# Row 43: T/U/V updates
$ws.Range("T43").Value = 262351.3154
$ws.Range("U43").Value = 1311756.5536
$ws.Range("V43").Value = 2657.6136

# Row 44: H column item description change
$ws.Range("H44").Value = "petroleum coke (graphitized)"

# Row 46: H column item description change
$ws.Range("H46").Value = "calcined petroleum coke"
